$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17 — shifts existing rows 17:77 down to 18:78
# and extends the used range to A1:T78 (matches dimension change in the diff).
$ws.Rows(17).Insert()

# Populate the newly inserted row 17 with the new data record.
$ws.Cells.Item(17, 1).Value = 10
$ws.Cells.Item(17, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(17, 3).Value = 'La Araucanía'
$ws.Cells.Item(17, 4).Value = 44998
$ws.Cells.Item(17, 5).Value = 9
$ws.Cells.Item(17, 6).Value = 'Fruta'
$ws.Cells.Item(17, 7).Value = 100107
$ws.Cells.Item(17, 8).Value = 'Otros'
$ws.Cells.Item(17, 9).Value = 100107011
$ws.Cells.Item(17, 10).Value = 'Tuna'
$ws.Cells.Item(17, 11).Value = 'Sin especificar'
$ws.Cells.Item(17, 12).Value = 'Primera'
$ws.Cells.Item(17, 13).Value = 35
$ws.Cells.Item(17, 14).Value = 15000
$ws.Cells.Item(17, 15).Value = 15000
$ws.Cells.Item(17, 16).Value = 15000
$ws.Cells.Item(17, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(17, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(17, 19).Value = 938
$ws.Cells.Item(17, 20).Value = 16
